$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$questions = @(
    'What is 1 + 1?',
    'What is 1 + 1?',
    'What is 1 + 1?',
    'What is 1 + 1?',
    'What is 1 + 1?',
    'What is 1 + 1?',
    'What is AI?',
    'What is 1 + 1?',
    'What is AI?'
)

$answers = @(
    'The answer to the question "What is 1 + 1?" is 2.',
    'The answer to the question "What is 1 + 1?" is 2.',
    '2',
    '2',
    'The answer to the question "What is 1 + 1?" is 2.',
    '2.',
    'I don''t have any information about "AI" in the provided documents. The GEO application seems to be focused on well log authoring, analysis, and reporting for petroleum geologists and engineers. If you could provide more context or clarify what you mean by "AI", I''ll do my best to help.',
    '2.',
    'I don''t have any information about "AI" in the provided documents. The GEO application seems to be focused on well log authoring, analysis, and reporting for petroleum geologists and engineers. If you could provide more context or clarify what you mean by "AI", I''ll do my best to help.'
)

$startRow = 13
for ($i = 0; $i -lt $questions.Length; $i++) {
    $r = $startRow + $i
    $cA = $ws.Cells.Item($r, 1)
    $cB = $ws.Cells.Item($r, 2)
    $cA.NumberFormat = "@"
    $cB.NumberFormat = "@"
    $cA.Value = $questions[$i]
    $cB.Value = $answers[$i]
    $cA.Style = "Normal"
    $cB.Style = "Normal"
}
